$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks")

# Two task rows (the "jnvjdnvds" row and the "ds c nc d" row) were removed
# from the Tasks list, leaving only the "q" task (previously the 4th row),
# which shifts up to become row 2.
$ws.Rows("2:3").Delete()
